# Empleados.xlsx - "No hice nada basicamente"
#
# - EMPLEADOS!D18 status changed from "A" to "I"
# - A new test record appended at row 43 (with a truly blank row 42 above it)
# - Selection left on D41 (last "real" data row) after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLEADOS")

# D18: ESTATUS A -> I
$ws.Range("D18").Value = "I"

# Leave row 42 present-but-blank (no cell data) between the existing data
# and the new test row, same as in the authored workbook.
$ws.Rows("42:42").OutlineLevel = 0

# New "Prueba" test record on row 43.
$ws.Range("A43").Value = 999999
$ws.Range("B43").Value = "Prueba"
$ws.Range("C43").Value = "Probado"
$ws.Range("D43").Value = "A"
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 0

# The new row's cells use the default (unstyled) format, not the centered
# style inherited from columns E/F.
$ws.Range("A43:F43").Style = "Normal"

# Final cursor/selection position.
$ws.Range("D41").Select()
